$d = $word.ActiveDocument

# Locate the paragraph that ends the "Requisitos" list (LOQ4095 ...),
# then remove the following paragraphs up through the "Ver no Jupiter ..."
# line and the "(c) 2020 ... Creative Commons Attribution" footer line
# (including the blank paragraph separating them from the requirement
# text), leaving the trailing blank paragraph and the page-break
# paragraph untouched.

$anchor = $null
$footer = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "LOQ4095") {
        $anchor = $p
    } elseif ($t -match "Powered by Jekyll") {
        $footer = $p
    }
}

$startPos = $anchor.Range.End
$endPos = $footer.Range.End

$r = $d.Range($startPos, $endPos)
$r.Delete()
